$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string allocation: set B2 first so the new
# "#1341342" text lands before "Euroservice" in the shared strings table,
# matching how the distributor name (typo fix) was the last value entered.
$ws.Range("B2").Value = "#1341342"
$ws.Range("A2").Value = "Euroservice"

# Explicitly (re-)apply the General number format to the edited cells so
# the format is stamped on them rather than left implicit.
$ws.Range("A2:B2").NumberFormat = "General"

# Move the active selection to A2 (where the edit happened).
$ws.Range("A2").Select() | Out-Null
